$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-10 Saturday" "2026-01-11 Sunday"

Replace-Text "28×72=" "32×11="
Replace-Text "95×19=" "99×26="
Replace-Text "64×19=" "63×14="
Replace-Text "42×76=" "17×22="
Replace-Text "38×11=" "19×86="
Replace-Text "39×54=" "91×69="
Replace-Text "77×81=" "56×77="
Replace-Text "72×82=" "14×17="
Replace-Text "17×49=" "67×51="
Replace-Text "66×61=" "49×39="
Replace-Text "34×15=" "55×76="
Replace-Text "34×54=" "72×72="
Replace-Text "43×71=" "16×19="
Replace-Text "60×95=" "93×74="
Replace-Text "75×43=" "30×59="
Replace-Text "81×43=" "24×85="
Replace-Text "55×64=" "44×92="
Replace-Text "79×88=" "16×82="
Replace-Text "55×71=" "63×74="
Replace-Text "54×84=" "30×19="
Replace-Text "26×88=" "81×37="
Replace-Text "65×22=" "63×74="
Replace-Text "21×54=" "65×41="
Replace-Text "32×24=" "73×31="
Replace-Text "66×12=" "60×40="

Write-Output "done"
